$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "HK_R_acc_G"
$ws.Range("A2").Value = 73.992673992674
$ws.Range("A3").Value = 73.504273504273513
$ws.Range("A4").Value = 73.544973544973544
$ws.Range("A5").Value = 69.108669108669105
$ws.Range("A6").Value = 71.794871794871796
$ws.Range("A7").Value = 73.056573056573058
$ws.Range("A8").Value = 73.748473748473756
$ws.Range("A9").Value = 74.277574277574274
$ws.Range("A10").Value = 73.341473341473346
$ws.Range("A11").Value = 73.015873015873012
$ws.Range("A12").Value = 74.155474155474153
$ws.Range("A13").Value = 76.923076923076934
$ws.Range("A14").Value = 72.730972730972738
$ws.Range("A15").Value = 73.504273504273513
$ws.Range("A16").Value = 72.771672771672769
$ws.Range("A17").Value = 75.091575091575095
$ws.Range("A18").Value = 76.353276353276357
$ws.Range("A19").Value = 76.109076109076113
$ws.Range("A20").Value = 75.213675213675216
$ws.Range("A21").Value = 72.527472527472526
$ws.Range("A22").Value = 75.295075295075293
$ws.Range("A23").Value = 75.986975986975992
$ws.Range("A24").Value = 75.661375661375658
$ws.Range("A25").Value = 75.905575905575901
$ws.Range("A26").Value = 73.789173789173788
$ws.Range("A27").Value = 73.789173789173788
$ws.Range("A28").Value = 74.074074074074076
$ws.Range("A29").Value = 80.179080179080174
$ws.Range("A30").Value = 79.527879527879534
$ws.Range("A31").Value = 77.574277574277573
$ws.Range("A32").Value = 71.550671550671552
$ws.Range("A33").Value = 74.033374033374031
$ws.Range("A34").Value = 74.033374033374031
$ws.Range("A35").Value = 75.539275539275536
$ws.Range("A36").Value = 75.661375661375658
$ws.Range("A37").Value = 79.242979242979246
$ws.Range("A38").Value = 75.579975579975581
$ws.Range("A39").Value = 75.45787545787546
$ws.Range("A40").Value = 76.271876271876266
$ws.Range("A41").Value = 74.643874643874639
$ws.Range("A42").Value = 74.358974358974365
$ws.Range("A43").Value = 74.765974765974761
$ws.Range("A44").Value = 73.951973951973955
$ws.Range("A45").Value = 74.236874236874243
$ws.Range("A46").Value = 73.300773300773301
$ws.Range("A47").Value = 72.85307285307286
$ws.Range("A48").Value = 72.486772486772495
$ws.Range("A49").Value = 73.015873015873012
